# Update cryptos price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.847.46'
$ws.Range('E2').Value = '  +0.15%  '
$ws.Range('D3').Value = '1.633.87'
$ws.Range('E3').Value = '  +0.33%  '
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.97'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.27%  '
$ws.Range('E6').Value = '  -0.31%  '
$ws.Range('E7').Value = '  -0.17%  '
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('E9').Value = '  -0.11%  '
$ws.Range('E10').Value = '  +2.87%  '
$ws.Range('E11').Value = '  +0.19%  '
$ws.Range('E12').Value = '  -0.07%  '
$ws.Range('D13').Value = '1.635.90'
$ws.Range('E13').Value = '  +0.36%  '
$ws.Range('D14').Value = '1.859.89'
$ws.Range('E14').Value = '  +0.32%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.559'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.46%  '
$ws.Range('D16').Value = '0.0₃0768'
$ws.Range('E16').Value = '  +2.12%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '62.98'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.46%  '
$ws.Range('D18').Value = '25.860.37'
$ws.Range('E18').Value = '  +0.15%  '
$ws.Range('E19').Value = '  -0.17%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '193.97'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.56%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.39'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.22%  '
$ws.Range('E22').Value = '  +1.23%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.22'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.76%  '
$ws.Range('E24').Value = '  -0.15%  '
$ws.Range('E25').Value = '  -3.77%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '139.04'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.76%  '
$ws.Range('E27').Value = '  -3.64%  '
$ws.Range('E28').Value = '  +1.60%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.52'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.78%  '
$ws.Range('E30').Value = '  +0.04%  '
$ws.Range('E31').Value = '  +2.04%  '
$ws.Range('E32').Value = '  +0.73%  '
$ws.Range('E33').Value = '  +1.98%  '
$ws.Range('E34').Value = '  +1.13%  '
$ws.Range('E35').Value = '  +0.22%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.902'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.02%  '
$ws.Range('E37').Value = '  +1.02%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.550'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.89%  '
$ws.Range('B39').Value = 'Maker'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D39').Value = '1.124.13'
$ws.Range('E39').Value = '  -0.88%  '
$ws.Range('E40').Value = '  +0.05%  '
$ws.Range('E41').Value = '  +0.77%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.51'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.00%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '99.66'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.56%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.801'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.84%  '
$ws.Range('E45').Value = '  -1.53%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '55.48'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.02%  '
$ws.Range('E47').Value = '  -5.02%  '
$ws.Range('B48').Value = 'SynthetixNetwork'
$ws.Range('C48').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.40'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +11.29%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.66'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.12%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0503'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.53%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.00'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.38%  '
